$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row 13 with the "Removed OptiFine Lightmaps" text, matching the
# existing formatting used by the other single-column rows (A4:A12).
$ws.Range("A13").Value = "Removed OptiFine Lightmaps"
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A13").Value = "Removed OptiFine Lightmaps"

# Move the active selection to A14, as reflected in the saved view state.
$ws.Range("A14").Select()
